$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 203 ("「小さな生き物にしては危険な名前だ」" / Tasmanian devil post),
# which shifts all subsequent rows up by one.
$ws.Rows.Item(203).Delete()
